$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A221").Value = 'I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted"'
$ws.Range("B221").Value = 'llama3.2:latest'
$ws.Range("C221").Value = @'
Here is an optimized answer:
I want to add this name as my curve shade name length "Hydrocarbon bearing zone highlighted".
To set the curve shading name, click on the "Curve" menu and select "Shading". Then, in the "Shading" dialog box, enter "Hydrocarbon bearing zone highlighted" in the "Name" field.
Note that the character limit for the curve shade name is 20 characters. Since your chosen name has a length of 37 characters, it exceeds the allowed limit. Therefore, this operation is not allowed due to the character length constraint.
This answer addresses the feedback by:
1. Providing a clear step-by-step guide on how to set the curve shading name.
2. Mentioning the character limit of 20 characters for the curve shade name and explaining why the operation is not allowed due to exceeding this limit.
3. Encouraging the user to check if their chosen name exceeds the character limit and make a decision accordingly.
This optimized answer should provide a more accurate and helpful response to the user's question while also addressing the concerns raised in the feedback.
'@

$ws.Range("A222").Value = 'How many tadpole definitions am I allowed?'
$ws.Range("B222").Value = 'llama3.2:latest'
$ws.Range("C222").Value = @'
Here is an optimized answer:
"I want to add this name as my curve shade name length \"Hydrocarbon bearing zone highlighted\".
To set the curve shading name, click on the \"Curve\" menu and select \"Shading\". Then, in the \"Shading\" dialog box, enter \"Hydrocarbon bearing zone highlighted\" in the \"Name\" field.
Note that the character limit for the curve shade name is 20 characters. Since your chosen name has a length of 37 characters, it exceeds the allowed limit. Therefore, this operation is not allowed due to the character length constraint."
'@

$ws.Range("A223").Value = 'How many tadpole definitions am I allowed?'
$ws.Range("B223").Value = 'llama3.2:latest'
$ws.Range("C223").Value = @'
Here is an optimized answer:
"I want to add this name as my curve shade name length \"Hydrocarbon bearing zone highlighted\".
To set the curve shading name, click on the \"Curve\" menu and select \"Shading\". Then, in the \"Shading\" dialog box, enter \"Hydrocarbon bearing zone highlighted\" in the \"Name\" field.
Note that the character limit for the curve shade name is 20 characters. Since your chosen name has a length of 37 characters, it exceeds the allowed limit. Therefore, this operation is not allowed due to the character length constraint."
'@

$ws.Range("A224").Value = 'How many tadpole definitions am I allowed?'
$ws.Range("B224").Value = 'llama3.2:latest'
$ws.Range("C224").Value = @'
According to the document, you are allowed to create a maximum of 50 new tadpole definition templates.
'@

# Undo the wrap-on-paste autofit so new rows keep the sheet default height,
# matching how the source rows (with embedded newlines) already look.
$ws.Range("221:224").EntireRow.AutoFit()
